# Automatic update of files.
#
# The underlying source export re-ordered a handful of observation rows.
# Re-apply the same row-content changes against the already-open workbook.
#
# Rows 6, 7, 8: a 3-cycle rotation of the per-record fields
#   (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn,
#    Auktor, Ost, Nord) -> new row6 = old row8, new row7 = old row6,
#    new row8 = old row7.
#
# Rows 37 <-> 38 and rows 39 <-> 40: simple swaps of the same fields
# (plus "Rodlistade" and the presence of the empty "Kon" cell for 37/38,
# which travel with the record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Record($row, $id, $sort, $taxid, $art, $vet, $auk, $ost, $nord) {
    $ws.Range("A$row").Value = $id
    $ws.Range("B$row").Value = $sort
    $ws.Range("E$row").Value = $taxid
    $ws.Range("F$row").Value = $art
    $ws.Range("G$row").Value = $vet
    $ws.Range("H$row").Value = $auk
    $ws.Range("Q$row").Value = $ost
    $ws.Range("R$row").Value = $nord
}

# ---- Rows 6 / 7 / 8 -------------------------------------------------------
Set-Record 6 111866919 90823 5966 "Motaggsvamp" "Sarcodon squamosus" "(Schaeff.) Quél." 703071 7299536
Set-Record 7 111867119 88623 1962 "Vaddporing" "Anomoporia kamtschatica" "(Parmasto) Bondartseva" 703180 7299491
Set-Record 8 111867113 90794 4362 "Blå taggsvamp" "Hydnellum caeruleum" "(Hornem.) P.Karst." 703089 7299467

# ---- Rows 37 / 38 (full swap, including Rodlistade + Kon placeholder) ----
$ws.Range("D37").Value = "NT"
$ws.Range("D38").Value = "LC"

Set-Record 37 111866994 90843 5448 "Svartvit taggsvamp" "Phellodon connatus" "(Schultz) nom.prov" 703115 7299511
Set-Record 38 111867271 95693 221941 "Plattlummer" "Lycopodium complanatum" "L." 703122 7299415

# The empty "Kon" (L) placeholder cell travels with the Plattlummer record.
$ws.Range("L37").ClearContents()
$ws.Range("L38").Value = ""

# ---- Rows 39 / 40 (full swap) ---------------------------------------------
Set-Record 39 111867456 90786 3100 "Talltaggsvamp" "Bankera fuligineoalba" "(Schmidt : Fr.) Pouzar" 703129 7299348
Set-Record 40 111867059 77388 6446 "Kolflarnlav" "Carbonicola anthracophila" "(Nyl.) Bendiksby & Timdal" 703135 7299505
